# Update the "想去人数" (number of people interested) counts on the
# 展览 (Exhibitions) and 全部类型 (All types) sheets, as generated by the
# gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 6399
$wsExpo.Range("F3").Value  = 112
$wsExpo.Range("F8").Value  = 531
$wsExpo.Range("F9").Value  = 77
$wsExpo.Range("F10").Value = 71
$wsExpo.Range("F13").Value = 365
$wsExpo.Range("F14").Value = 775
$wsExpo.Range("F15").Value = 3116
$wsExpo.Range("F17").Value = 183
$wsExpo.Range("F18").Value = 1779

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 6399
$wsAll.Range("F3").Value  = 112
$wsAll.Range("F9").Value  = 531
$wsAll.Range("F10").Value = 77
$wsAll.Range("F11").Value = 71
$wsAll.Range("F14").Value = 365
$wsAll.Range("F15").Value = 775
$wsAll.Range("F16").Value = 3116
$wsAll.Range("F18").Value = 183
$wsAll.Range("F19").Value = 1779
